$d = $word.ActiveDocument

# --- Locate the insertion point: right before the closing ")" that follows
# "...city)" in the first (User) bullet paragraph. ---
$searchRange = $d.Content
$found = $searchRange.Find.Execute("city)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'city)' in document"
}
$closeParenPos = $searchRange.End - 1

# Insert the new field text just before the ")".
$insertRange = $d.Range($closeParenPos, $closeParenPos)
$insertRange.InsertBefore(", incomeid (fk)")

# --- Split the newly-inserted text into its own run, distinct from the
# preceding "...city" run, by briefly bookmarking the boundary and then
# removing that temporary bookmark (this forces a run split at that
# position without merging formatting/rsid metadata across the
# boundary). ---
$boundaryPos = $closeParenPos
$boundaryRange = $d.Range($boundaryPos, $boundaryPos)
[void]$d.Bookmarks.Add("_zzTmpSplit", $boundaryRange)
$tmpBm = $d.Bookmarks.Item("_zzTmpSplit")
$tmpBm.Delete()

# --- Re-add the "_GoBack" bookmark at the true new last-edit location,
# i.e. immediately before the final ")" (now shifted by the inserted
# text). Word keeps only one "_GoBack" bookmark document-wide, so adding
# it here implicitly removes/relocates any pre-existing "_GoBack"
# elsewhere (e.g. the one that used to sit before "profit" in the
# Income bullet). ---
$newClose = $boundaryPos + (", incomeid (fk)").Length
$goBackRange = $d.Range($newClose, $newClose)
[void]$d.Bookmarks.Add("_GoBack", $goBackRange)
